$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column numeric-looking text (e.g. "314.41") must be forced to Text via NumberFormat="@"
# before assignment, then ClearFormats() restores the default (unstyled) cell format so the
# only observable change is the cell's value, matching the source data's inline-string typing.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.898.23'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.238.26'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.65%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.41'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.61'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.01%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -7.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.99'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0821'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.50%  '

$ws.Range("E12").Value = '  -7.20%  '

$ws.Range("E13").Value = '  -2.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.578.79'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.241.19'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.839'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.88'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -5.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.770.48'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -7.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0973'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.75'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.71'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.49%  '

$ws.Range("E24").Value = '  -7.32%  '

$ws.Range("E25").Value = '  -8.34%  '

$ws.Range("E26").Value = '  +0.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("E28").Value = '  -2.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.46'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.96'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -8.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.40'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0830'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -6.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.33'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.65'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.91'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -8.23%  '

$ws.Range("E37").Value = '  -5.69%  '

$ws.Range("E38").Value = '  -3.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.55'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -12.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -10.78%  '

$ws.Range("E42").Value = '  -6.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.701.67'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.37'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.07%  '

$ws.Range("E46").Value = '  -6.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.15'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.60'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '71.40'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '56.20'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.61'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.42%  '
